$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 08:22"

# Israel (row 25) - updated case counts
$ws.Range("B25").Value = 13883
$ws.Range("C25").Value = 170
$ws.Range("D25").Value = 4353
$ws.Range("E25").Value = 9349
$ws.Range("F25").Value = 142
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 181

# Oman (row 73) - updated case counts
$ws.Range("E73").Value = 1164
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 8

# Bulgaria (row 85) - updated case counts
$ws.Range("B85").Value = 966
$ws.Range("C85").Value = 37
$ws.Range("D85").Value = 170
$ws.Range("E85").Value = 751
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 45

# Reunion / Georgia swapped order + Georgia's data updated, Reunion keeps its
# previous figures (just moves down one row)
$ws.Range("A109").Value = "Georgia"
$ws.Range("B109").Value = 408
$ws.Range("C109").Value = 6
$ws.Range("D109").Value = 95
$ws.Range("E109").Value = 309
$ws.Range("F109").Value = 6
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 4

$ws.Range("A110").Value = "Reunion"
$ws.Range("B110").Value = 408
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 238
$ws.Range("E110").Value = 170
$ws.Range("F110").Value = 2
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 0

# Martinica (row 129) - updated case counts
$ws.Range("E129").Value = 76
$ws.Range("G129").Value = 2
$ws.Range("H129").Value = 14
